$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.973.12'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '2.603.36'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.37'
$ws.Range("E5").Value = '  +3.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.89'
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  +2.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.70'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.348'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("D13").Value = '3.055.33'
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").Value = '60.992.59'
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.74'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '2.604.68'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.75'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '353.96'
$ws.Range("E19").Value = '  +2.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.59'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.23'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.16'
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.428'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.716.83'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.994'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").Value = '0.0₃0849'
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.33'
$ws.Range("E31").Value = '  +10.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.39'
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.61'
$ws.Range("E33").Value = '  +3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.28'
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.24'
$ws.Range("E35").Value = '  +6.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.949'
$ws.Range("E36").Value = '  +10.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.21'
$ws.Range("E37").Value = '  +2.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.50'
$ws.Range("E38").Value = '  +2.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.849'
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.80'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.47'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '287.82'
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("E43").Value = '  +1.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.626'
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0561'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.00'
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.57'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0238'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.03'
$ws.Range("E51").Value = '  +7.95%  '
